$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.360.95'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '3.520.93'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '569.77'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '182.53'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.18%  '
$ws.Range('D7').Value = '3.517.74'
$ws.Range('E7').Value = '  -1.49%  '
$ws.Range('E8').Value = '  -2.70%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.186'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.94%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.639'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.99%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '53.56'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -5.12%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000299'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.45'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.24%  '
$ws.Range('D15').Value = '4.090.90'
$ws.Range('E15').Value = '  -1.32%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '19.25'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.20%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '69.264.88'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.510.86'
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.33'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '541.08'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +14.76%  '
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '19.52'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.93'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.56%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.35'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '93.68'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +5.87%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.13'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.90'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.00%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.09'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.85%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '31.69'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.36'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -4.67%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '12.52'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.60%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '65.04'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('E34').Value = '  -5.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '572.17'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.12'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +6.32%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '38.01'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.97%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.399'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').Value = '0.0₃0763'
$ws.Range('E40').Value = '  -4.95%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.11'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.134'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -6.77%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.37'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.98%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.49'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.84%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.96'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.40%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0443'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').Value = '3.138.05'
$ws.Range('E47').Value = '  -2.71%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.18'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -4.04%  '
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('E51').Value = '  +19.94%  '
